# Apply updated cryptos data (prices and volume %) to sheet1, per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.264.31"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.606.26"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.80"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.42"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "1.829.33"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "1.608.06"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "26.230.10"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.07"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "200.38"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.99"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.20"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").Value = "1.166.01"
$ws.Range("E36").Value = "  +4.27%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.35"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "1.740.57"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.14"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0106"
$ws.Range("E46").Value = "  +14.94%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.54"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.11"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("E51").Value = "  -0.15%  "
